$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy formatting/style from an existing header cell (A1) to the new header cells
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in team record (Wins/Losses/Ties) for every data row (2-43)
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 84
    $ws.Cells.Item($r, 31).Value = 78
    $ws.Cells.Item($r, 32).Value = 0
}
